$wb = $excel.ActiveWorkbook

# --- 1. Add the "Note" column to the CustomModsTable (Master Schedule sheet),
#        positioned between "Week" and "Time". ---
$ws5 = $wb.Worksheets.Item("Master Schedule")
$lo = $ws5.ListObjects.Item(1)

# This runtime's ListColumns.Add always appends at the end of the table,
# so append a column and then re-label the header cells (D/E/F) into the
# desired final order. Writing through HeaderRowRange keeps both the sheet
# cell and the table's column name in sync.
$lo.ListColumns.Add(4) | Out-Null
$lo.HeaderRowRange.Cells.Item(1, 4).Value = "Note"
$lo.HeaderRowRange.Cells.Item(1, 5).Value = "Time"
$lo.HeaderRowRange.Cells.Item(1, 6).Value = "Equipment"

# --- 2. Selection bookkeeping that accompanies the fix: the user had been
#        working on Equipment (selecting B29) and finished on Master
#        Schedule (selecting I30), which is now the active tab. ---
$ws3 = $wb.Worksheets.Item("Inventory")
$ws3.Activate()
$ws3.Range("I3").Select()

$ws4 = $wb.Worksheets.Item("Equipment")
$ws4.Activate()
$ws4.Range("B29").Select()

$ws5.Activate()
$ws5.Range("I30").Select()
